# Change in unit of AIC and run of new results
#
# Each per-year worksheet ("2000".."2100") holds an AIC-style results block
# in rows 5:8, columns D:G. The commit rescales every one of those non-zero
# result values by a factor of 1e-6 (i.e. divide by 1,000,000) to reflect a
# change in unit, leaving untouched cells (labels, zeros, headers) as-is.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    for ($r = 5; $r -le 8; $r++) {
        for ($c = 4; $c -le 7; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            $val = $cell.Value2
            if ($val -ne 0) {
                $cell.Value = $val / 1000000
            }
        }
    }
}
